$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.277687430381775
$ws.Range("B1").Value = 2.319268941879272
$ws.Range("D1").Value = 1.386290311813354
$ws.Range("E1").Value = 0.8471335172653198
